$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.837.28"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "3.153.11"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'216.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.22%  "

$ws.Range("D6").Value = "'617.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.14%  "

$ws.Range("D7").Value = "'0.286"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +16.13%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.32%  "

$ws.Range("D10").Value = "3.148.87"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "'0.596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("D12").Value = "'0.0000254"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").Value = "'5.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.71%  "

$ws.Range("D15").Value = "3.735.68"
$ws.Range("E15").Value = "  -1.52%  "

$ws.Range("D16").Value = "'32.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Value = "81.780.07"
$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").Value = "3.148.74"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "'3.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.41%  "

$ws.Range("D20").Value = "'13.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.21%  "

$ws.Range("D21").Value = "'433.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("D22").Value = "'8.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.49%  "

$ws.Range("D23").Value = "'5.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.69%  "

$ws.Range("D24").Value = "'7.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").Value = "'5.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.40%  "

$ws.Range("D26").Value = "'11.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.09%  "

$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").Value = "'76.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("D30").Value = "'0.0000121"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'8.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.03%  "

$ws.Range("D33").Value = "'565.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.32%  "

$ws.Range("D34").Value = "'1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("E35").Value = "  +18.15%  "

$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.60%  "

$ws.Range("D38").Value = "'22.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.79%  "

$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'6.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.52%  "

$ws.Range("D41").Value = "'0.404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("D42").Value = "'20.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.95%  "

$ws.Range("D43").Value = "'3.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.80%  "

$ws.Range("D44").Value = "'2.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.16%  "

$ws.Range("D45").Value = "'158.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "'186.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("D48").Value = "'44.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("D49").Value = "'1.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("D50").Value = "'26.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "

$ws.Range("E51").Value = "  -5.76%  "
